$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(594, 1).Value2 = "2022-03-25 13:47:48"
$ws.Cells.Item(594, 2).Value2 = "8: 255`n"
$ws.Cells.Item(594, 3).Value2 = 3600
$ws.Cells.Item(594, 4).Value2 = 1648212470
$ws.Cells.Item(594, 5).Value2 = 48
$ws.Cells.Item(594, 6).Value2 = "10F872226797"
$ws.Cells.Item(594, 7).Value2 = 175
$ws.Cells.Item(594, 8).Value2 = 12
$ws.Rows.Item(594).AutoFit()

$ws.Cells.Item(595, 1).Value2 = "2022-03-25 14:00:50"
$ws.Cells.Item(595, 2).Value2 = "8: 255`n"
$ws.Cells.Item(595, 3).Value2 = 3600
$ws.Cells.Item(595, 4).Value2 = 1648213252
$ws.Cells.Item(595, 5).Value2 = 54
$ws.Cells.Item(595, 6).Value2 = "10F872226797"
$ws.Cells.Item(595, 7).Value2 = 179
$ws.Cells.Item(595, 8).Value2 = 15
$ws.Rows.Item(595).AutoFit()

$ws.Cells.Item(596, 1).Value2 = "2022-03-25 14:04:52"
$ws.Cells.Item(596, 2).Value2 = "8: 255`n"
$ws.Cells.Item(596, 3).Value2 = 3600
$ws.Cells.Item(596, 4).Value2 = 1648213494
$ws.Cells.Item(596, 5).Value2 = 50
$ws.Cells.Item(596, 6).Value2 = "10F872226797"
$ws.Cells.Item(596, 7).Value2 = 187
$ws.Cells.Item(596, 8).Value2 = 18
$ws.Rows.Item(596).AutoFit()

$ws.Cells.Item(597, 1).Value2 = "2022-03-25 14:13:55"
$ws.Cells.Item(597, 2).Value2 = "8: 255`n"
$ws.Cells.Item(597, 3).Value2 = 3600
$ws.Cells.Item(597, 4).Value2 = 1648214037
$ws.Cells.Item(597, 5).Value2 = 44
$ws.Cells.Item(597, 6).Value2 = "10F872226797"
$ws.Cells.Item(597, 7).Value2 = 177
$ws.Cells.Item(597, 8).Value2 = 21
$ws.Rows.Item(597).AutoFit()

$ws.Cells.Item(598, 1).Value2 = "2022-06-09 13:48:07"
$ws.Cells.Item(598, 2).Value2 = "8: 255`n"
$ws.Cells.Item(598, 3).Value2 = 3600
$ws.Cells.Item(598, 4).Value2 = 1654775301
$ws.Cells.Item(598, 5).Value2 = 50
$ws.Cells.Item(598, 6).Value2 = "10F872226797"
$ws.Cells.Item(598, 7).Value2 = 180
$ws.Cells.Item(598, 8).Value2 = 11
$ws.Rows.Item(598).AutoFit()

$ws.Cells.Item(599, 1).Value2 = "2022-06-09 16:48:28"
$ws.Cells.Item(599, 2).Value2 = "8: 255`n"
$ws.Cells.Item(599, 3).Value2 = 3600
$ws.Cells.Item(599, 4).Value2 = 1654786123
$ws.Cells.Item(599, 5).Value2 = 46
$ws.Cells.Item(599, 6).Value2 = "10F872226797"
$ws.Cells.Item(599, 7).Value2 = 188
$ws.Cells.Item(599, 8).Value2 = 28
$ws.Rows.Item(599).AutoFit()

$ws.Cells.Item(600, 1).Value2 = "2022-06-09 19:06:44"
$ws.Cells.Item(600, 2).Value2 = "8: 255`n"
$ws.Cells.Item(600, 3).Value2 = 3600
$ws.Cells.Item(600, 4).Value2 = 1654794418
$ws.Cells.Item(600, 5).Value2 = 46
$ws.Cells.Item(600, 6).Value2 = "10F872226797"
$ws.Cells.Item(600, 7).Value2 = 217
$ws.Cells.Item(600, 8).Value2 = 36
$ws.Rows.Item(600).AutoFit()

$ws.Cells.Item(601, 1).Value2 = "2022-06-09 22:55:57"
$ws.Cells.Item(601, 2).Value2 = "8: 255`n"
$ws.Cells.Item(601, 3).Value2 = 3600
$ws.Cells.Item(601, 4).Value2 = 1654808172
$ws.Cells.Item(601, 5).Value2 = 46
$ws.Cells.Item(601, 6).Value2 = "10F872226797"
$ws.Cells.Item(601, 7).Value2 = -57255
$ws.Cells.Item(601, 8).Value2 = 44
$ws.Rows.Item(601).AutoFit()

$ws.Cells.Item(602, 1).Value2 = "2022-06-09 22:56:17"
$ws.Cells.Item(602, 2).Value2 = "8: 255`n"
$ws.Cells.Item(602, 3).Value2 = 3600
$ws.Cells.Item(602, 4).Value2 = 1654808193
$ws.Cells.Item(602, 5).Value2 = 46
$ws.Cells.Item(602, 6).Value2 = "10F872226797"
$ws.Cells.Item(602, 7).Value2 = 222
$ws.Cells.Item(602, 8).Value2 = 44
$ws.Rows.Item(602).AutoFit()

$ws.Cells.Item(603, 1).Value2 = "2022-06-09 22:57:20"
$ws.Cells.Item(603, 2).Value2 = "8: 255`n"
$ws.Cells.Item(603, 3).Value2 = 3600
$ws.Cells.Item(603, 4).Value2 = 1654808255
$ws.Cells.Item(603, 5).Value2 = 46
$ws.Cells.Item(603, 6).Value2 = "10F872226797"
$ws.Cells.Item(603, 7).Value2 = 231
$ws.Cells.Item(603, 8).Value2 = 44
$ws.Rows.Item(603).AutoFit()

$ws.Cells.Item(604, 1).Value2 = "2022-06-09 22:58:17"
$ws.Cells.Item(604, 2).Value2 = "8: 255`n"
$ws.Cells.Item(604, 3).Value2 = 3600
$ws.Cells.Item(604, 4).Value2 = 1654808313
$ws.Cells.Item(604, 5).Value2 = 46
$ws.Cells.Item(604, 6).Value2 = "10F872226797"
$ws.Cells.Item(604, 7).Value2 = 263
$ws.Cells.Item(604, 8).Value2 = 44
$ws.Rows.Item(604).AutoFit()

$ws.Cells.Item(605, 1).Value2 = "2022-06-09 22:59:20"
$ws.Cells.Item(605, 2).Value2 = "8: 255`n"
$ws.Cells.Item(605, 3).Value2 = 3600
$ws.Cells.Item(605, 4).Value2 = 1654808375
$ws.Cells.Item(605, 5).Value2 = 46
$ws.Cells.Item(605, 6).Value2 = "10F872226797"
$ws.Cells.Item(605, 7).Value2 = 206
$ws.Cells.Item(605, 8).Value2 = 44
$ws.Rows.Item(605).AutoFit()

$ws.Cells.Item(606, 1).Value2 = "2022-06-10 09:18:11"
$ws.Cells.Item(606, 2).Value2 = "8: 255`n"
$ws.Cells.Item(606, 3).Value2 = 3600
$ws.Cells.Item(606, 4).Value2 = 1654845491
$ws.Cells.Item(606, 5).Value2 = 60
$ws.Cells.Item(606, 6).Value2 = "10F872226797"
$ws.Cells.Item(606, 7).Value2 = 182
$ws.Cells.Item(606, 8).Value2 = 4
$ws.Rows.Item(606).AutoFit()

$ws.Cells.Item(607, 1).Value2 = "2022-06-10 10:14:38"
$ws.Cells.Item(607, 2).Value2 = "8: 255`n"
$ws.Cells.Item(607, 3).Value2 = 3600
$ws.Cells.Item(607, 4).Value2 = 1654848878
$ws.Cells.Item(607, 5).Value2 = 60
$ws.Cells.Item(607, 6).Value2 = "10F872226797"
$ws.Cells.Item(607, 7).Value2 = 211
$ws.Cells.Item(607, 8).Value2 = 4
$ws.Rows.Item(607).AutoFit()

$ws.Cells.Item(608, 1).Value2 = "2022-06-10 13:52:02"
$ws.Cells.Item(608, 2).Value2 = "8: 255`n"
$ws.Cells.Item(608, 3).Value2 = 3600
$ws.Cells.Item(608, 4).Value2 = 1654861924
$ws.Cells.Item(608, 5).Value2 = 52
$ws.Cells.Item(608, 6).Value2 = "10F872226797"
$ws.Cells.Item(608, 7).Value2 = 189
$ws.Cells.Item(608, 8).Value2 = 3
$ws.Rows.Item(608).AutoFit()

$ws.Cells.Item(609, 1).Value2 = "2022-06-10 14:48:29"
$ws.Cells.Item(609, 2).Value2 = "8: 255`n"
$ws.Cells.Item(609, 3).Value2 = 3600
$ws.Cells.Item(609, 4).Value2 = 1654865310
$ws.Cells.Item(609, 5).Value2 = 52
$ws.Cells.Item(609, 6).Value2 = "10F872226797"
$ws.Cells.Item(609, 7).Value2 = 262
$ws.Cells.Item(609, 8).Value2 = 4
$ws.Rows.Item(609).AutoFit()

$ws.Cells.Item(610, 1).Value2 = "2022-06-10 15:04:30"
$ws.Cells.Item(610, 2).Value2 = "8: 255`n"
$ws.Cells.Item(610, 3).Value2 = 3600
$ws.Cells.Item(610, 4).Value2 = 1654866272
$ws.Cells.Item(610, 5).Value2 = 52
$ws.Cells.Item(610, 6).Value2 = "10F872226797"
$ws.Cells.Item(610, 7).Value2 = 192
$ws.Cells.Item(610, 8).Value2 = 4
$ws.Rows.Item(610).AutoFit()

$ws.Cells.Item(611, 1).Value2 = "2022-06-10 16:30:16"
$ws.Cells.Item(611, 2).Value2 = "8: 255`n"
$ws.Cells.Item(611, 3).Value2 = 3600
$ws.Cells.Item(611, 4).Value2 = 1654871418
$ws.Cells.Item(611, 5).Value2 = 58
$ws.Cells.Item(611, 6).Value2 = "10F872226797"
$ws.Cells.Item(611, 7).Value2 = 204
$ws.Cells.Item(611, 8).Value2 = 14
$ws.Rows.Item(611).AutoFit()

$ws.Cells.Item(612, 1).Value2 = "2022-06-10 18:50:26"
$ws.Cells.Item(612, 2).Value2 = "8: 255`n"
$ws.Cells.Item(612, 3).Value2 = 3600
$ws.Cells.Item(612, 4).Value2 = 1654879828
$ws.Cells.Item(612, 5).Value2 = 44
$ws.Cells.Item(612, 6).Value2 = "10F872226797"
$ws.Cells.Item(612, 7).Value2 = 190
$ws.Cells.Item(612, 8).Value2 = 22
$ws.Rows.Item(612).AutoFit()

$ws.Cells.Item(613, 1).Value2 = "2022-06-10 18:56:26"
$ws.Cells.Item(613, 2).Value2 = "8: 255`n"
$ws.Cells.Item(613, 3).Value2 = 3600
$ws.Cells.Item(613, 4).Value2 = 1654880189
$ws.Cells.Item(613, 5).Value2 = 44
$ws.Cells.Item(613, 6).Value2 = "10F872226797"
$ws.Cells.Item(613, 7).Value2 = 256
$ws.Cells.Item(613, 8).Value2 = 25
$ws.Rows.Item(613).AutoFit()

$ws.Cells.Item(614, 1).Value2 = "2022-06-10 18:56:56"
$ws.Cells.Item(614, 2).Value2 = "8: 255`n"
$ws.Cells.Item(614, 3).Value2 = 3600
$ws.Cells.Item(614, 4).Value2 = 1654880189
$ws.Cells.Item(614, 5).Value2 = 44
$ws.Cells.Item(614, 6).Value2 = "10F872226797"
$ws.Cells.Item(614, 7).Value2 = 256
$ws.Cells.Item(614, 8).Value2 = 25
$ws.Rows.Item(614).AutoFit()

$ws.Cells.Item(615, 1).Value2 = "2022-06-10 18:57:26"
$ws.Cells.Item(615, 2).Value2 = "8: 255`n"
$ws.Cells.Item(615, 3).Value2 = 3600
$ws.Cells.Item(615, 4).Value2 = 1654880189
$ws.Cells.Item(615, 5).Value2 = 44
$ws.Cells.Item(615, 6).Value2 = "10F872226797"
$ws.Cells.Item(615, 7).Value2 = 256
$ws.Cells.Item(615, 8).Value2 = 25
$ws.Rows.Item(615).AutoFit()

$ws.Cells.Item(616, 1).Value2 = "2022-06-10 18:58:06"
$ws.Cells.Item(616, 2).Value2 = "8: 255`n"
$ws.Cells.Item(616, 3).Value2 = 3600
$ws.Cells.Item(616, 4).Value2 = 1654880289
$ws.Cells.Item(616, 5).Value2 = 46
$ws.Cells.Item(616, 6).Value2 = "10F872226797"
$ws.Cells.Item(616, 7).Value2 = -360292
$ws.Cells.Item(616, 8).Value2 = 28
$ws.Rows.Item(616).AutoFit()

$ws.Cells.Item(617, 1).Value2 = "2022-06-10 18:58:36"
$ws.Cells.Item(617, 2).Value2 = "8: 255`n"
$ws.Cells.Item(617, 3).Value2 = 3600
$ws.Cells.Item(617, 4).Value2 = 1654880289
$ws.Cells.Item(617, 5).Value2 = 46
$ws.Cells.Item(617, 6).Value2 = "10F872226797"
$ws.Cells.Item(617, 7).Value2 = -360292
$ws.Cells.Item(617, 8).Value2 = 28
$ws.Rows.Item(617).AutoFit()

$ws.Cells.Item(618, 1).Value2 = "2022-06-11 08:02:11"
$ws.Cells.Item(618, 2).Value2 = "8: 255`n"
$ws.Cells.Item(618, 3).Value2 = 3600
$ws.Cells.Item(618, 4).Value2 = 1654927337
$ws.Cells.Item(618, 5).Value2 = 52
$ws.Cells.Item(618, 6).Value2 = "10F872226797"
$ws.Cells.Item(618, 7).Value2 = 192
$ws.Cells.Item(618, 8).Value2 = 7
$ws.Rows.Item(618).AutoFit()

$ws.Cells.Item(619, 1).Value2 = "2022-06-11 08:35:35"
$ws.Cells.Item(619, 2).Value2 = "8: 255`n"
$ws.Cells.Item(619, 3).Value2 = 3600
$ws.Cells.Item(619, 4).Value2 = 1654929340
$ws.Cells.Item(619, 5).Value2 = 52
$ws.Cells.Item(619, 6).Value2 = "10F872226797"
$ws.Cells.Item(619, 7).Value2 = 227
$ws.Cells.Item(619, 8).Value2 = 7
$ws.Rows.Item(619).AutoFit()

$ws.Cells.Item(620, 1).Value2 = "2022-06-11 09:07:37"
$ws.Cells.Item(620, 2).Value2 = "8: 255`n"
$ws.Cells.Item(620, 3).Value2 = 3600
$ws.Cells.Item(620, 4).Value2 = 1654931263
$ws.Cells.Item(620, 5).Value2 = 52
$ws.Cells.Item(620, 6).Value2 = "10F872226797"
$ws.Cells.Item(620, 7).Value2 = 190
$ws.Cells.Item(620, 8).Value2 = 7
$ws.Rows.Item(620).AutoFit()

$ws.Cells.Item(621, 1).Value2 = "2022-06-24 12:59:38"
$ws.Cells.Item(621, 2).Value2 = "8: 255`n"
$ws.Cells.Item(621, 3).Value2 = 3600
$ws.Cells.Item(621, 4).Value2 = 1656068379
$ws.Cells.Item(621, 5).Value2 = 60
$ws.Cells.Item(621, 6).Value2 = "10F872226797"
$ws.Cells.Item(621, 7).Value2 = 285
$ws.Cells.Item(621, 8).Value2 = 6
$ws.Rows.Item(621).AutoFit()

$ws.Cells.Item(622, 1).Value2 = "2022-06-24 13:06:01"
$ws.Cells.Item(622, 2).Value2 = "8: 255`n"
$ws.Cells.Item(622, 3).Value2 = 3600
$ws.Cells.Item(622, 4).Value2 = 1656068761
$ws.Cells.Item(622, 5).Value2 = 60
$ws.Cells.Item(622, 6).Value2 = "10F872226797"
$ws.Cells.Item(622, 7).Value2 = 208
$ws.Cells.Item(622, 8).Value2 = 6
$ws.Rows.Item(622).AutoFit()

$ws.Cells.Item(623, 1).Value2 = "2022-06-24 18:39:29"
$ws.Cells.Item(623, 2).Value2 = "8: 255`n"
$ws.Cells.Item(623, 3).Value2 = 3600
$ws.Cells.Item(623, 4).Value2 = 1656088771
$ws.Cells.Item(623, 5).Value2 = 60
$ws.Cells.Item(623, 6).Value2 = "10F872226797"
$ws.Cells.Item(623, 7).Value2 = 184
$ws.Cells.Item(623, 8).Value2 = 11
$ws.Rows.Item(623).AutoFit()

$ws.Cells.Item(624, 1).Value2 = "2022-06-25 00:36:26"
$ws.Cells.Item(624, 2).Value2 = "8: 255`n"
$ws.Cells.Item(624, 3).Value2 = 3600
$ws.Cells.Item(624, 4).Value2 = 1656110189
$ws.Cells.Item(624, 5).Value2 = 60
$ws.Cells.Item(624, 6).Value2 = "10F872226797"
$ws.Cells.Item(624, 7).Value2 = -59843
$ws.Cells.Item(624, 8).Value2 = 14
$ws.Rows.Item(624).AutoFit()

$ws.Cells.Item(625, 1).Value2 = "2022-06-25 00:36:46"
$ws.Cells.Item(625, 2).Value2 = "8: 255`n"
$ws.Cells.Item(625, 3).Value2 = 3600
$ws.Cells.Item(625, 4).Value2 = 1656110210
$ws.Cells.Item(625, 5).Value2 = 60
$ws.Cells.Item(625, 6).Value2 = "10F872226797"
$ws.Cells.Item(625, 7).Value2 = 181
$ws.Cells.Item(625, 8).Value2 = 14
$ws.Rows.Item(625).AutoFit()

$ws.Cells.Item(626, 1).Value2 = "2022-06-25 00:37:46"
$ws.Cells.Item(626, 2).Value2 = "8: 255`n"
$ws.Cells.Item(626, 3).Value2 = 3600
$ws.Cells.Item(626, 4).Value2 = 1656110270
$ws.Cells.Item(626, 5).Value2 = 60
$ws.Cells.Item(626, 6).Value2 = "10F872226797"
$ws.Cells.Item(626, 7).Value2 = 185
$ws.Cells.Item(626, 8).Value2 = 14
$ws.Rows.Item(626).AutoFit()

$ws.Cells.Item(627, 1).Value2 = "2022-06-25 00:38:46"
$ws.Cells.Item(627, 2).Value2 = "8: 255`n"
$ws.Cells.Item(627, 3).Value2 = 3600
$ws.Cells.Item(627, 4).Value2 = 1656110330
$ws.Cells.Item(627, 5).Value2 = 60
$ws.Cells.Item(627, 6).Value2 = "10F872226797"
$ws.Cells.Item(627, 7).Value2 = 194
$ws.Cells.Item(627, 8).Value2 = 14
$ws.Rows.Item(627).AutoFit()

$ws.Cells.Item(628, 1).Value2 = "2022-06-25 06:32:07"
$ws.Cells.Item(628, 2).Value2 = "8: 255`n"
$ws.Cells.Item(628, 3).Value2 = 3600
$ws.Cells.Item(628, 4).Value2 = 1656131526
$ws.Cells.Item(628, 5).Value2 = 58
$ws.Cells.Item(628, 6).Value2 = "10F872226797"
$ws.Cells.Item(628, 7).Value2 = 176
$ws.Cells.Item(628, 8).Value2 = 6
$ws.Rows.Item(628).AutoFit()

$ws.Cells.Item(629, 1).Value2 = "2022-06-25 06:40:29"
$ws.Cells.Item(629, 2).Value2 = "8: 255`n"
$ws.Cells.Item(629, 3).Value2 = 3600
$ws.Cells.Item(629, 4).Value2 = 1656132028
$ws.Cells.Item(629, 5).Value2 = 58
$ws.Cells.Item(629, 6).Value2 = "10F872226797"
$ws.Cells.Item(629, 7).Value2 = 177
$ws.Cells.Item(629, 8).Value2 = 6
$ws.Rows.Item(629).AutoFit()

$ws.Cells.Item(630, 1).Value2 = "2022-06-25 06:40:59"
$ws.Cells.Item(630, 2).Value2 = "8: 255`n"
$ws.Cells.Item(630, 3).Value2 = 3600
$ws.Cells.Item(630, 4).Value2 = 1656132028
$ws.Cells.Item(630, 5).Value2 = 58
$ws.Cells.Item(630, 6).Value2 = "10F872226797"
$ws.Cells.Item(630, 7).Value2 = 177
$ws.Cells.Item(630, 8).Value2 = 6
$ws.Rows.Item(630).AutoFit()

$ws.Cells.Item(631, 1).Value2 = "2022-06-25 11:08:34"
$ws.Cells.Item(631, 2).Value2 = "8: 255`n"
$ws.Cells.Item(631, 3).Value2 = 3600
$ws.Cells.Item(631, 4).Value2 = 1656148113
$ws.Cells.Item(631, 5).Value2 = 60
$ws.Cells.Item(631, 6).Value2 = "10F872226797"
$ws.Cells.Item(631, 7).Value2 = 167
$ws.Cells.Item(631, 8).Value2 = 10
$ws.Rows.Item(631).AutoFit()

$ws.Cells.Item(632, 1).Value2 = "2022-06-25 11:26:57"
$ws.Cells.Item(632, 2).Value2 = "8: 255`n"
$ws.Cells.Item(632, 3).Value2 = 3600
$ws.Cells.Item(632, 4).Value2 = 1656149216
$ws.Cells.Item(632, 5).Value2 = 60
$ws.Cells.Item(632, 6).Value2 = "10F872226797"
$ws.Cells.Item(632, 7).Value2 = 178
$ws.Cells.Item(632, 8).Value2 = 10
$ws.Rows.Item(632).AutoFit()

$ws.Cells.Item(633, 1).Value2 = "2022-06-25 11:54:00"
$ws.Cells.Item(633, 2).Value2 = "8: 255`n"
$ws.Cells.Item(633, 3).Value2 = 3600
$ws.Cells.Item(633, 4).Value2 = 1656150839
$ws.Cells.Item(633, 5).Value2 = 60
$ws.Cells.Item(633, 6).Value2 = "10F872226797"
$ws.Cells.Item(633, 7).Value2 = 178
$ws.Cells.Item(633, 8).Value2 = 10
$ws.Rows.Item(633).AutoFit()

$ws.Cells.Item(634, 1).Value2 = "2022-06-25 12:55:06"
$ws.Cells.Item(634, 2).Value2 = "8: 255`n"
$ws.Cells.Item(634, 3).Value2 = 3600
$ws.Cells.Item(634, 4).Value2 = 1656154506
$ws.Cells.Item(634, 5).Value2 = 60
$ws.Cells.Item(634, 6).Value2 = "10F872226797"
$ws.Cells.Item(634, 7).Value2 = 180
$ws.Cells.Item(634, 8).Value2 = 10
$ws.Rows.Item(634).AutoFit()

$ws.Cells.Item(635, 1).Value2 = "2022-07-01 16:36:24"
$ws.Cells.Item(635, 2).Value2 = "8: 255`n"
$ws.Cells.Item(635, 3).Value2 = 3600
$ws.Cells.Item(635, 4).Value2 = 1656686192
$ws.Cells.Item(635, 5).Value2 = 60
$ws.Cells.Item(635, 6).Value2 = "10F872226797"
$ws.Cells.Item(635, 7).Value2 = 196
$ws.Cells.Item(635, 8).Value2 = 16
$ws.Rows.Item(635).AutoFit()

$ws.Cells.Item(636, 1).Value2 = "2022-07-01 17:16:31"
$ws.Cells.Item(636, 2).Value2 = "8: 255`n"
$ws.Cells.Item(636, 3).Value2 = 3600
$ws.Cells.Item(636, 4).Value2 = 1656688600
$ws.Cells.Item(636, 5).Value2 = 52
$ws.Cells.Item(636, 6).Value2 = "10F872226797"
$ws.Cells.Item(636, 7).Value2 = 180
$ws.Cells.Item(636, 8).Value2 = 28
$ws.Rows.Item(636).AutoFit()

$ws.Cells.Item(637, 1).Value2 = "2022-07-01 17:40:34"
$ws.Cells.Item(637, 2).Value2 = "8: 255`n"
$ws.Cells.Item(637, 3).Value2 = 3600
$ws.Cells.Item(637, 4).Value2 = 1656690043
$ws.Cells.Item(637, 5).Value2 = 62
$ws.Cells.Item(637, 6).Value2 = "10F872226797"
$ws.Cells.Item(637, 7).Value2 = 936
$ws.Cells.Item(637, 8).Value2 = 31
$ws.Rows.Item(637).AutoFit()

Write-Output "done"